$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 441.86957
$ws.Range("I53").Value = 241.5
$ws.Range("J53").Value = 899.8570999999999
$ws.Range("K53").Value = 241.5
$ws.Range("L53").Value = 899.8570999999999
$ws.Range("M53").Value = 395.5
$ws.Range("N53").Value = -2173.8571
$ws.Range("H98").Value = 7685.75
$ws.Range("I98").Value = 6021.1113
$ws.Range("J98").Value = 8474.263000000001
$ws.Range("K98").Value = 6021.1113
$ws.Range("L98").Value = 8474.263000000001
$ws.Range("M98").Value = -4523.1113
$ws.Range("N98").Value = -11470.263
$ws.Range("H106").Value = 2485.35
$ws.Range("I106").Value = 1994.6471
$ws.Range("K106").Value = 1994.6471
$ws.Range("M106").Value = -1363.6471
$ws.Range("H122").Value = 7685.75
$ws.Range("I122").Value = 6021.1113
$ws.Range("J122").Value = 8474.263000000001
$ws.Range("K122").Value = 18063.3339
$ws.Range("L122").Value = 25422.789
$ws.Range("M122").Value = -15613.3339
$ws.Range("N122").Value = -30322.789
$ws.Range("H138").Value = 3061.5574
$ws.Range("I138").Value = 1560.1818
$ws.Range("J138").Value = 3908.487
$ws.Range("K138").Value = 4680.5454
$ws.Range("L138").Value = 11725.461
$ws.Range("M138").Value = 459.4546
$ws.Range("N138").Value = -22005.461
$ws.Range("H141").Value = 6788.227
$ws.Range("I141").Value = 7072.737
$ws.Range("J141").Value = 4986.3335
$ws.Range("K141").Value = 21218.211
$ws.Range("L141").Value = 14959.0005
$ws.Range("M141").Value = -16038.211
$ws.Range("N141").Value = -25319.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3306
$ws.Range("I32").Value = 3135.0735
$ws.Range("J32").Value = 3859.476
$ws.Range("K32").Value = 3135.0735
$ws.Range("L32").Value = 3859.476
$ws.Range("M32").Value = -2848.0735
$ws.Range("N32").Value = -4433.476000000001
$ws.Range("H102").Value = 1775
$ws.Range("I102").Value = 1775
$ws.Range("K102").Value = 1775
$ws.Range("M102").Value = -153
$ws.Range("H137").Value = 44458.332
$ws.Range("J137").Value = 44458.332
$ws.Range("L137").Value = 44458.332
$ws.Range("N137").Value = -54658.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H141").Value = 39889
$ws.Range("J141").Value = 39889
$ws.Range("L141").Value = 39889
$ws.Range("N141").Value = -50249

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2184.0657
$ws.Range("I31").Value = 942.5
$ws.Range("J31").Value = 3106.3713
$ws.Range("K31").Value = 942.5
$ws.Range("L31").Value = 3106.3713
$ws.Range("M31").Value = -647.5
$ws.Range("N31").Value = -3696.3713
$ws.Range("H34").Value = 2184.0657
$ws.Range("I34").Value = 942.5
$ws.Range("J34").Value = 3106.3713
$ws.Range("K34").Value = 942.5
$ws.Range("L34").Value = 3106.3713
$ws.Range("M34").Value = -740.5
$ws.Range("N34").Value = -3510.3713
$ws.Range("H124").Value = 47264.715
$ws.Range("J124").Value = 47264.715
$ws.Range("L124").Value = 47264.715
$ws.Range("N124").Value = -52174.715
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H127").Value = 43250
$ws.Range("J127").Value = 43250
$ws.Range("L127").Value = 43250
$ws.Range("N127").Value = -53170
$ws.Range("H134").Value = 1466.5186
$ws.Range("I134").Value = 959
$ws.Range("J134").Value = 5526.6665
$ws.Range("K134").Value = 2877
$ws.Range("L134").Value = 16579.9995
$ws.Range("M134").Value = -342
$ws.Range("N134").Value = -21649.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1587.6227
$ws.Range("J68").Value = 1665.6586
$ws.Range("L68").Value = 4996.9758
$ws.Range("N68").Value = -6618.9758
$ws.Range("H71").Value = 1587.6227
$ws.Range("J71").Value = 1665.6586
$ws.Range("L71").Value = 14990.9274
$ws.Range("N71").Value = -23102.9274

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 23887.6
$ws.Range("J46").Value = 24928
$ws.Range("L46").Value = 24928
$ws.Range("N46").Value = -25240
$ws.Range("H70").Value = 6841.0884
$ws.Range("I70").Value = 6087.52
$ws.Range("J70").Value = 8934.333000000001
$ws.Range("K70").Value = 6087.52
$ws.Range("L70").Value = 8934.333000000001
$ws.Range("M70").Value = -5817.52
$ws.Range("N70").Value = -9474.333000000001
$ws.Range("H73").Value = 6841.0884
$ws.Range("I73").Value = 6087.52
$ws.Range("J73").Value = 8934.333000000001
$ws.Range("K73").Value = 6087.52
$ws.Range("L73").Value = 8934.333000000001
$ws.Range("M73").Value = -5151.52
$ws.Range("N73").Value = -10806.333
$ws.Range("H80").Value = 12502999
$ws.Range("I80").Value = 19233876
$ws.Range("J80").Value = 2800.8572
$ws.Range("K80").Value = 19233876
$ws.Range("L80").Value = 2800.8572
$ws.Range("M80").Value = -19232878
$ws.Range("N80").Value = -4796.8572
$ws.Range("H83").Value = 12502999
$ws.Range("I83").Value = 19233876
$ws.Range("J83").Value = 2800.8572
$ws.Range("K83").Value = 96169380
$ws.Range("L83").Value = 14004.286
$ws.Range("M83").Value = -96164388
$ws.Range("N83").Value = -23988.286
$ws.Range("H132").Value = 3875.7778
$ws.Range("I132").Value = 2211.4614
$ws.Range("J132").Value = 5421.2144
$ws.Range("K132").Value = 6634.3842
$ws.Range("L132").Value = 16263.6432
$ws.Range("M132").Value = -4104.3842
$ws.Range("N132").Value = -21323.6432
$ws.Range("H134").Value = 32713.928
$ws.Range("J134").Value = 32713.928
$ws.Range("L134").Value = 98141.784
$ws.Range("N134").Value = -103211.784
$ws.Range("H135").Value = 47707.06
$ws.Range("J135").Value = 47707.06
$ws.Range("L135").Value = 47707.06
$ws.Range("N135").Value = -57847.06
$ws.Range("H137").Value = 29500
$ws.Range("J137").Value = 43750
$ws.Range("L137").Value = 43750
$ws.Range("N137").Value = -53950

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3132.8572
$ws.Range("I132").Value = 2357.842
$ws.Range("J132").Value = 5810.1816
$ws.Range("K132").Value = 7073.526
$ws.Range("L132").Value = 17430.5448
$ws.Range("M132").Value = -4543.526
$ws.Range("N132").Value = -22490.5448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4553.375
$ws.Range("I122").Value = 2647.25
$ws.Range("K122").Value = 7941.75
$ws.Range("M122").Value = -5491.75
$ws.Range("H132").Value = 8131978.5
$ws.Range("I132").Value = 1199.4828
$ws.Range("J132").Value = 27781362
$ws.Range("K132").Value = 3598.4484
$ws.Range("L132").Value = 83344086
$ws.Range("M132").Value = -1068.4484
$ws.Range("N132").Value = -83349146
